$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: Replace the italic "meta description" text that currently
# sits in the very last paragraph of the document with the new
# image-prompt text. We assign straight to the sub-range's .Text
# (rather than using Find.Execute's replacement, which would run the
# text through AutoCorrect/"smart quotes" and mangle the apostrophes)
# so the run keeps its existing italic formatting and the apostrophes
# stay as plain ASCII characters.
# ------------------------------------------------------------------
$newImageText = "Create a feature image for Dragon Spin Pick n Mix that features a happy Maya warrior with glasses in a cartoon-style. The warrior should be surrounded by dragons and treasure, creating a sense of adventure and excitement. The image should be bright and colorful, with an Eastern-inspired design to match the game's theme. The Maya warrior should be depicted as confident and triumphant, holding a winning jackpot symbol above their head. This will appeal to the slot game's adventurous and adventurous players, while also highlighting the exciting gameplay and potential for big wins."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$imgRange = $lastPara.Range
$imgRange.End = $imgRange.End - 1
$imgRange.Text = $newImageText

# ------------------------------------------------------------------
# Step 2: Remove the duplicated bold heading paragraph
# ("Play Dragon Spin Pick n Mix Free - Review of Features & Payouts")
# that was left sitting right before the final (now updated)
# paragraph -- this duplicate is no longer needed since the title
# text will live in the new "Meta description" paragraph at the top.
# We skip paragraph 1 (the real Heading1 title) and only remove the
# later plain-paragraph duplicate.
# ------------------------------------------------------------------
$titleText = "Play Dragon Spin Pick n Mix Free - Review of Features & Payouts"
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 2; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq ($titleText + [char]13)) {
        $para.Range.Delete()
    }
}

# ------------------------------------------------------------------
# Step 3: Insert a new "Meta description" paragraph right after the
# document's first paragraph (the Heading1 title). It contains a
# bold "Meta description" run followed by a plain run with the
# description text, matching Word's native OOXML shape (with the
# leading empty run that normal body paragraphs in this document
# carry).
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Read a review of Dragon Spin Pick n Mix, an online slot game with varying bet sizes and high volatility. Play for free and activate multiple bonuses.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml) | Out-Null
